# Updated cryptos list - price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "42.692.37"
Set-TextCell $ws.Range("E2") "  +1.09%  "

Set-TextCell $ws.Range("D3") "2.300.36"
Set-TextCell $ws.Range("E3") "  +0.00%  "

Set-TextCell $ws.Range("E4") "  +0.16%  "

Set-TextCell $ws.Range("D5") "316.03"
Set-TextCell $ws.Range("E5") "  -0.43%  "

Set-TextCell $ws.Range("D6") "103.05"
Set-TextCell $ws.Range("E6") "  -0.11%  "

Set-TextCell $ws.Range("E7") "  -0.33%  "

Set-TextCell $ws.Range("E8") "  +0.23%  "

Set-TextCell $ws.Range("E9") "  -0.76%  "

Set-TextCell $ws.Range("D10") "39.62"
Set-TextCell $ws.Range("E10") "  -0.12%  "

Set-TextCell $ws.Range("D11") "0.0908"
Set-TextCell $ws.Range("E11") "  -0.05%  "

Set-TextCell $ws.Range("D12") "8.55"
Set-TextCell $ws.Range("E12") "  +2.59%  "

Set-TextCell $ws.Range("E14") "  +4.32%  "

Set-TextCell $ws.Range("D15") "15.36"
Set-TextCell $ws.Range("E15") "  +0.56%  "

Set-TextCell $ws.Range("D16") "2.643.84"
Set-TextCell $ws.Range("E16") "  -0.20%  "

Set-TextCell $ws.Range("D17") "2.294.85"
Set-TextCell $ws.Range("E17") "  +0.50%  "

Set-TextCell $ws.Range("D18") "42.589.25"
Set-TextCell $ws.Range("E18") "  +0.82%  "

Set-TextCell $ws.Range("E19") "  +1.75%  "

Set-TextCell $ws.Range("E20") "  +0.35%  "

Set-TextCell $ws.Range("D21") "13.77"
Set-TextCell $ws.Range("E21") "  +21.71%  "

Set-TextCell $ws.Range("D22") "74.01"
Set-TextCell $ws.Range("E22") "  +0.58%  "

Set-TextCell $ws.Range("E23") "  -3.18%  "

Set-TextCell $ws.Range("D24") "267.21"
Set-TextCell $ws.Range("E24") "  -4.78%  "

Set-TextCell $ws.Range("D25") "2.24"
Set-TextCell $ws.Range("E25") "  -1.16%  "

Set-TextCell $ws.Range("D26") "1.00"
Set-TextCell $ws.Range("E26") "  +0.09%  "

Set-TextCell $ws.Range("D27") "10.93"
Set-TextCell $ws.Range("E27") "  +0.36%  "

Set-TextCell $ws.Range("E28") "  -3.52%  "

Set-TextCell $ws.Range("D29") "22.65"
Set-TextCell $ws.Range("E29") "  -1.54%  "

Set-TextCell $ws.Range("D30") "6.63"
Set-TextCell $ws.Range("E30") "  +12.49%  "

Set-TextCell $ws.Range("D31") "37.19"
Set-TextCell $ws.Range("E31") "  +2.53%  "

Set-TextCell $ws.Range("D32") "165.56"
Set-TextCell $ws.Range("E32") "  +0.69%  "

Set-TextCell $ws.Range("D33") "0.0883"
Set-TextCell $ws.Range("E33") "  +1.00%  "

Set-TextCell $ws.Range("D34") "0.132"
Set-TextCell $ws.Range("E34") "  -3.07%  "

Set-TextCell $ws.Range("D35") "2.56"
Set-TextCell $ws.Range("E35") "  -3.16%  "

Set-TextCell $ws.Range("E36") "  -1.58%  "

Set-TextCell $ws.Range("D37") "4.56"
Set-TextCell $ws.Range("E37") "  -1.05%  "

Set-TextCell $ws.Range("E38") "  +1.43%  "

Set-TextCell $ws.Range("D39") "3.75"
Set-TextCell $ws.Range("E39") "  -0.08%  "

Set-TextCell $ws.Range("E41") "  +9.54%  "

Set-TextCell $ws.Range("D42") "70.68"
Set-TextCell $ws.Range("E42") "  +2.00%  "

Set-TextCell $ws.Range("D43") "96.14"
Set-TextCell $ws.Range("E43") "  -3.46%  "

Set-TextCell $ws.Range("E45") "  +0.39%  "

Set-TextCell $ws.Range("D46") "12.46"
Set-TextCell $ws.Range("E46") "  +3.85%  "

Set-TextCell $ws.Range("D47") "117.33"
Set-TextCell $ws.Range("E47") "  +4.77%  "

Set-TextCell $ws.Range("D48") "80.22"
Set-TextCell $ws.Range("E48") "  +3.64%  "

Set-TextCell $ws.Range("D49") "1.663.39"
Set-TextCell $ws.Range("E49") "  +3.89%  "

Set-TextCell $ws.Range("D50") "5.29"
Set-TextCell $ws.Range("E50") "  +0.00%  "

Set-TextCell $ws.Range("D51") "8.86"
Set-TextCell $ws.Range("E51") "  -0.55%  "
